# Generate Report for Handback
# Adds the newly handed-back file "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md"
# as a new row (row 4) to the Overview, zh-cn and de-de tables.

$wb = $excel.ActiveWorkbook

# Cornflowerblue (FF6495ED) expressed as an OLE BGR color for Font.Color,
# matching the workbook's existing custom "HyperLink" cell style.
$hyperlinkColor = 15570276
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

function Style-AsDate($range) {
    $range.NumberFormat = $dateFormat
}

# ---------------------------------------------------------------------------
# Sheet "Overview" - add row 4
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md"
$wsOverview.Range("B4").Value = "e2e\fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-10-10 09:30:12"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fcee92e4e213c4854bdcd3f29e3f80fcfe213c4/e2e/fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md", "", "", "e2e\fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md") | Out-Null
Style-AsHyperlink $wsOverview.Range("B4")
Style-AsDate $wsOverview.Range("G4")

# ---------------------------------------------------------------------------
# Sheet "zh-cn" - add row 4
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A4").Value = "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "True"
$wsZhCn.Range("G4").Value = "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.1e9488531f109ac96f0c2db36c9a1a28b29c3380.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-10-10 09:30:00"
$wsZhCn.Range("I4").Value = "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md"
$wsZhCn.Range("J4").Value = "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.1e9488531f109ac96f0c2db36c9a1a28b29c3380.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-10-10 09:30:42"
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("O4").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fcee92e4e213c4854bdcd3f29e3f80fcfe213c4/e2e/fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md", "", "", "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fcee92e4e213c4854bdcd3f29e3f80fcfe213c4/e2e/fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md", "", "", "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md") | Out-Null
Style-AsHyperlink $wsZhCn.Range("A4")
Style-AsHyperlink $wsZhCn.Range("I4")
Style-AsDate $wsZhCn.Range("H4")
Style-AsDate $wsZhCn.Range("K4")

# ---------------------------------------------------------------------------
# Sheet "de-de" - add row 4
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A4").Value = "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "True"
$wsDeDe.Range("G4").Value = "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.1e9488531f109ac96f0c2db36c9a1a28b29c3380.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-10-10 09:30:12"
$wsDeDe.Range("I4").Value = "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md"
$wsDeDe.Range("J4").Value = "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.1e9488531f109ac96f0c2db36c9a1a28b29c3380.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-10-10 09:30:58"
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("O4").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fcee92e4e213c4854bdcd3f29e3f80fcfe213c4/e2e/fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md", "", "", "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fcee92e4e213c4854bdcd3f29e3f80fcfe213c4/e2e/fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md", "", "", "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md") | Out-Null
Style-AsHyperlink $wsDeDe.Range("A4")
Style-AsHyperlink $wsDeDe.Range("I4")
Style-AsDate $wsDeDe.Range("H4")
Style-AsDate $wsDeDe.Range("K4")

Write-Output "Added handback row for fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md to Overview, zh-cn and de-de sheets."
